# AChangeSettings.pptx - "Add files via upload" edit replay
#
# 1) Notes Master "Date Placeholder" - datetimeFigureOut field text
#    12/3/2019 -> 7/6/2021 (presentation was re-saved on 2021-07-06).
# 2) Slide 4 ("General Settings"), Content Placeholder 2, bullet
#    "Version - The version of the translationStudio program is displayed."
#    -> "Version - The version of the BTT Writer program is displayed."
#    (re-branding translationStudio -> BTT Writer), which PowerPoint stores
#    by splitting the run that used to hold the whole sentence.

$p = $ppt.ActivePresentation

# --- 1) Notes master date placeholder -------------------------------------
$notesMaster = $p.NotesMaster
$dateAndTime = $notesMaster.HeadersFooters.DateAndTime
$dateAndTime.Text = "7/6/2021"

# --- 2) Slide 4 bullet text -------------------------------------------------
$slide = $p.Slides.Item(4)
$shape = $slide.Shapes.Item(3)
$bodyRange = $shape.TextFrame.TextRange
$paras = $bodyRange.Paragraphs()

# The bullet currently reads:
#   "Version" + " - The version of the translationStudio program is displayed."
# Replace just the product name, which splits the trailing run into three
# runs: "... the ", "BTT Writer program ", "is displayed."
$oldWord = "translationStudio program "
$newWord = "BTT Writer program "

for ($i = 1; $i -le $paras.Count; $i++) {
    $para = $bodyRange.Paragraphs($i, 1)
    $paraText = $para.Text
    if ($paraText.StartsWith("Version") -and $paraText.Contains($oldWord)) {
        $startIdx = $paraText.IndexOf($oldWord)
        $target = $para.Characters($startIdx + 1, $oldWord.Length)
        $target.Text = $newWord
        break
    }
}
